$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-46 for columns B and C with new values
$ws.Cells.Item(2, 2).Value = 0.3057755019566411
$ws.Cells.Item(2, 3).Value = 0.9826765718905971
$ws.Cells.Item(3, 2).Value = 3.154138662649057
$ws.Cells.Item(3, 3).Value = 2.116306662292992
$ws.Cells.Item(4, 2).Value = 6.309303283115023
$ws.Cells.Item(4, 3).Value = 3.212001523098469
$ws.Cells.Item(5, 2).Value = 6.429515363513585
$ws.Cells.Item(5, 3).Value = 4.135790728667687
$ws.Cells.Item(6, 2).Value = 7.639755830693706
$ws.Cells.Item(6, 3).Value = 5.023303323608117
$ws.Cells.Item(7, 2).Value = 9.813262399729423
$ws.Cells.Item(7, 3).Value = 6.147740111681724
$ws.Cells.Item(8, 2).Value = 10.07208983726656
$ws.Cells.Item(8, 3).Value = 7.223257217080174
$ws.Cells.Item(9, 2).Value = 12.81391831371356
$ws.Cells.Item(9, 3).Value = 8.614052548334671
$ws.Cells.Item(10, 2).Value = 13.54667688665792
$ws.Cells.Item(10, 3).Value = 9.641044288093667
$ws.Cells.Item(11, 2).Value = 18.70980162941933
$ws.Cells.Item(11, 3).Value = 10.75538690925083
$ws.Cells.Item(12, 2).Value = 18.84208360291822
$ws.Cells.Item(12, 3).Value = 11.74825053739646
$ws.Cells.Item(13, 2).Value = 20.29284793574165
$ws.Cells.Item(13, 3).Value = 12.88617047708607
$ws.Cells.Item(14, 2).Value = 21.99655093544908
$ws.Cells.Item(14, 3).Value = 14.0047038236492
$ws.Cells.Item(15, 2).Value = 23.47394361526146
$ws.Cells.Item(15, 3).Value = 15.12988592880661
$ws.Cells.Item(16, 2).Value = 28.61399442879771
$ws.Cells.Item(16, 3).Value = 16.44029817557212
$ws.Cells.Item(17, 2).Value = 30.34999258236576
$ws.Cells.Item(17, 3).Value = 17.49082622723459
$ws.Cells.Item(18, 2).Value = 32.6333819200124
$ws.Cells.Item(18, 3).Value = 18.49820391976089
$ws.Cells.Item(19, 2).Value = 32.68631422510158
$ws.Cells.Item(19, 3).Value = 19.48321643640386
$ws.Cells.Item(20, 2).Value = 40.20928467167986
$ws.Cells.Item(20, 3).Value = 20.55627083938649
$ws.Cells.Item(21, 2).Value = 40.27763790002835
$ws.Cells.Item(21, 3).Value = 21.78894017037864
$ws.Cells.Item(22, 2).Value = 40.49940116259118
$ws.Cells.Item(22, 3).Value = 23.28908568810289
$ws.Cells.Item(23, 2).Value = 40.75847444187139
$ws.Cells.Item(23, 3).Value = 24.34671023788813
$ws.Cells.Item(24, 2).Value = 43.94925535551054
$ws.Cells.Item(24, 3).Value = 25.6788196804849
$ws.Cells.Item(25, 2).Value = 44.45507674979105
$ws.Cells.Item(25, 3).Value = 26.96685689275536
$ws.Cells.Item(26, 2).Value = 46.9090320987316
$ws.Cells.Item(26, 3).Value = 27.94268979019827
$ws.Cells.Item(27, 2).Value = 47.55470329123317
$ws.Cells.Item(27, 3).Value = 29.05747652407612
$ws.Cells.Item(28, 2).Value = 47.99606648034335
$ws.Cells.Item(28, 3).Value = 30.04198406806526
$ws.Cells.Item(29, 2).Value = 48.07506795306357
$ws.Cells.Item(29, 3).Value = 31.20628923462507
$ws.Cells.Item(30, 2).Value = 51.38584688565206
$ws.Cells.Item(30, 3).Value = 32.33574363433733
$ws.Cells.Item(31, 2).Value = 52.19765464963719
$ws.Cells.Item(31, 3).Value = 33.46645357182354
$ws.Cells.Item(32, 2).Value = 54.33709891283433
$ws.Cells.Item(32, 3).Value = 34.47790722735392
$ws.Cells.Item(33, 2).Value = 55.05536286757207
$ws.Cells.Item(33, 3).Value = 35.6349914708116
$ws.Cells.Item(34, 2).Value = 55.69204492397582
$ws.Cells.Item(34, 3).Value = 36.79903387750633
$ws.Cells.Item(35, 2).Value = 60.55470327552327
$ws.Cells.Item(35, 3).Value = 37.93367368211315
$ws.Cells.Item(36, 2).Value = 62.31901260921842
$ws.Cells.Item(36, 3).Value = 38.91201586919249
$ws.Cells.Item(37, 2).Value = 62.44778901826707
$ws.Cells.Item(37, 3).Value = 40.09100077052401
$ws.Cells.Item(38, 2).Value = 62.80807289820739
$ws.Cells.Item(38, 3).Value = 41.18956802005567
$ws.Cells.Item(39, 2).Value = 76.53776748414062
$ws.Cells.Item(39, 3).Value = 42.24328014358694
$ws.Cells.Item(40, 2).Value = 76.64724414875722
$ws.Cells.Item(40, 3).Value = 43.30882812755814
$ws.Cells.Item(41, 2).Value = 78.7631379500087
$ws.Cells.Item(41, 3).Value = 44.36228994056093
$ws.Cells.Item(42, 2).Value = 80.27799027442124
$ws.Cells.Item(42, 3).Value = 45.33509096520394
$ws.Cells.Item(43, 2).Value = 83.72280812048349
$ws.Cells.Item(43, 3).Value = 46.56263662006835
$ws.Cells.Item(44, 2).Value = 83.95777005763956
$ws.Cells.Item(44, 3).Value = 47.56359617577904
$ws.Cells.Item(45, 2).Value = 86.98944131623007
$ws.Cells.Item(45, 3).Value = 49.14259820764222
$ws.Cells.Item(46, 2).Value = 92.14298144941921
$ws.Cells.Item(46, 3).Value = 50.37186127014811

# Add new rows 47-50 (A values 45-48) with column A styled like the existing column A cells
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 92.21644588337426
$ws.Cells.Item(47, 3).Value = 51.53160675000726
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 92.78042424572141
$ws.Cells.Item(48, 3).Value = 52.72856272289629
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 94.21593395708028
$ws.Cells.Item(49, 3).Value = 53.82659609278578
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 94.52484625553295
$ws.Cells.Item(50, 3).Value = 54.97474088191304

# Copy style from an existing column-A cell (e.g. A46) down into the new A47:A50 cells
$ws.Range("A46").Copy()
$ws.Range("A47:A50").PasteSpecial(-4122)
$excel.CutCopyMode = 0
